$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 12692
$ws.Range("B2").Value = "Nicolas Cardoso"
$ws.Range("C2").Value = "TI"
$ws.Range("D2").Value = "Consulta medica"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45104
$ws.Range("G2").Value = 4867.85

# Row 3
$ws.Range("A3").Value = 31024
$ws.Range("B3").Value = "Maria Júlia Macedo"
$ws.Range("C3").Value = "TI"
$ws.Range("D3").Value = "Consulta medica"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 45080
$ws.Range("G3").Value = 5158.75

# Row 4
$ws.Range("A4").Value = 28693
$ws.Range("B4").Value = "Erick Marques"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Viagem de negocios"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45091
$ws.Range("G4").Value = 3668.5

# Row 5
$ws.Range("A5").Value = 12015
$ws.Range("B5").Value = "Bianca Casa Grande"
$ws.Range("C5").Value = "P&D"
$ws.Range("D5").Value = "Viagem de negocios"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45085
$ws.Range("G5").Value = 8695.07

# Row 6
$ws.Range("A6").Value = 82415
$ws.Range("B6").Value = "Mariane Freitas"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Viagem de negocios"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45102
$ws.Range("G6").Value = 7443.73

# Row 7
$ws.Range("A7").Value = 1461
$ws.Range("B7").Value = "Maria Clara Pereira"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45088
$ws.Range("G7").Value = 8983.83

# Row 8
$ws.Range("A8").Value = 68126
$ws.Range("B8").Value = "Natália Alves"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("D8").Value = "Consulta medica"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45104
$ws.Range("G8").Value = 9026.459999999999

# Row 9
$ws.Range("A9").Value = 17216
$ws.Range("B9").Value = "Davi Miguel Borges"
$ws.Range("C9").Value = "P&D"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45081
$ws.Range("G9").Value = 2161.78

# Row 10
$ws.Range("A10").Value = 20948
$ws.Range("B10").Value = "Luiz Henrique da Rosa"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Consulta medica"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45085
$ws.Range("G10").Value = 5059.31

# Row 11
$ws.Range("A11").Value = 5527
$ws.Range("B11").Value = "Gustavo Henrique da Mata"
$ws.Range("C11").Value = "Juridico"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45099
$ws.Range("G11").Value = 7791.76
